$wb = $excel.ActiveWorkbook

# Rename the "Task Name" column header to "Name" on the PV-Test-01 sheet.
$ws1 = $wb.Worksheets.Item("PV-Test-01")

$ws1.Range("C1").Value = "Name"

# Make PV-Test-01 the active sheet/tab (was "Dummy"), and update its
# selected cell from C4 to C2.
$ws1.Activate()
$ws1.Range("C2").Select()
